# correccion en df de autos y motos
# Updates occupancy counts (OCUPACION AUTOS / OCUPACION MOTOS columns D-I)
# across all five parking-lot sheets to match the corrected source data.

$wb = $excel.ActiveWorkbook

# --- CENTENARIO ---
$ws = $wb.Worksheets.Item("CENTENARIO")
$ws.Range("G10").Value = 25
$ws.Range("G11").Value = 1
$ws.Range("G17").Value = 3

# --- EL_PEÑON ---
$ws = $wb.Worksheets.Item("EL_PEÑON")
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 2
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("G14").Value = 46
$ws.Range("F15").Value = 19
$ws.Range("I15").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("G18").Value = 8
$ws.Range("D24").Value = 23
$ws.Range("F30").Value = 19
$ws.Range("I30").Value = 1
$ws.Range("D32").Value = 30
$ws.Range("G32").Value = 11
$ws.Range("D40").Value = 47
$ws.Range("E40").Value = 53
$ws.Range("G40").Value = 43
$ws.Range("H40").Value = 1
$ws.Range("D43").Value = 2
$ws.Range("G43").Value = 9
$ws.Range("E46").Value = 27
$ws.Range("H46").Value = 1
$ws.Range("D53").Value = 38

# --- GRANADA ---
$ws = $wb.Worksheets.Item("GRANADA")
$ws.Range("D10").Value = 38
$ws.Range("G10").Value = 4
$ws.Range("I16").Value = 23
$ws.Range("D17").Value = 269
$ws.Range("G17").Value = 1
$ws.Range("D23").Value = 19
$ws.Range("G23").Value = 15
$ws.Range("I25").Value = 8
$ws.Range("D29").Value = 29
$ws.Range("G29").Value = 9
$ws.Range("H29").Value = 6
$ws.Range("H30").Value = 6
$ws.Range("D31").Value = 41
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = 7
$ws.Range("D32").Value = 21
$ws.Range("G32").Value = 8
$ws.Range("D38").Value = 68
$ws.Range("G38").Value = 1
$ws.Range("I40").Value = 23
$ws.Range("G42").Value = 8
$ws.Range("D47").Value = 34
$ws.Range("G47").Value = 15
$ws.Range("G52").Value = 19
$ws.Range("D68").Value = 92
$ws.Range("H68").Value = 9
$ws.Range("I68").Value = 8
$ws.Range("D69").Value = 40
$ws.Range("G69").Value = 30
$ws.Range("D73").Value = 52
$ws.Range("D75").Value = 11
$ws.Range("G75").Value = 16
$ws.Range("I75").Value = 22
$ws.Range("I78").Value = 10
$ws.Range("D79").Value = 20
$ws.Range("G79").Value = 1
$ws.Range("D81").Value = 38
$ws.Range("G81").Value = 1
$ws.Range("D86").Value = 5
$ws.Range("G86").Value = 23
$ws.Range("I89").Value = 93
$ws.Range("G91").Value = 10
$ws.Range("D92").Value = 15

# --- SAN_ANTONIO ---
$ws = $wb.Worksheets.Item("SAN_ANTONIO")
$ws.Range("D7").Value = 23
$ws.Range("G7").Value = 10
$ws.Range("G20").Value = 16
$ws.Range("D30").Value = 8
$ws.Range("D58").Value = 13
$ws.Range("G58").Value = 11
$ws.Range("D59").Value = 29
$ws.Range("G59").Value = 12
$ws.Range("D60").Value = 9
$ws.Range("D78").Value = 3
$ws.Range("G78").Value = 8

# --- SAN_FERNANDO_PARQUE_DEL_PERRO ---
$ws = $wb.Worksheets.Item("SAN_FERNANDO_PARQUE_DEL_PERRO")
$ws.Range("D6").Value = 20
$ws.Range("G6").Value = 12
$ws.Range("E8").Value = 16
$ws.Range("H8").Value = 1
$ws.Range("G18").Value = 6
$ws.Range("E19").Value = 2
$ws.Range("H19").Value = 1
$ws.Range("D22").Value = 46
$ws.Range("G22").Value = 5
$ws.Range("G28").Value = 10
$ws.Range("E31").Value = 7
$ws.Range("H31").Value = 1
$ws.Range("F35").Value = 10
$ws.Range("I35").Value = 1
$ws.Range("D36").Value = 15
$ws.Range("G36").Value = 1
$ws.Range("D46").Value = 9
$ws.Range("G46").Value = 11
$ws.Range("D48").Value = 20
$ws.Range("E48").Value = 11
$ws.Range("H48").Value = 1
$ws.Range("G55").Value = 6
$ws.Range("G60").Value = 11
$ws.Range("E65").Value = 21
$ws.Range("H65").Value = 1
$ws.Range("D67").Value = 46
$ws.Range("G67").Value = 5
$ws.Range("E72").Value = 28
$ws.Range("H72").Value = 1
$ws.Range("G77").Value = 0
$ws.Range("D78").Value = 22
$ws.Range("G78").Value = 4
$ws.Range("D93").Value = 18
$ws.Range("D94").Value = 7
$ws.Range("G94").Value = 19
$ws.Range("D99").Value = 24
$ws.Range("G99").Value = 1

